# update taakverdeling & logboek
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New task row (row 16): description, date, estimated time, actual time
$ws.Range("D16").Value = "Commentaar, properties en testen bij Tankaart"
$ws.Range("E16").Value = Get-Date -Year 2021 -Month 10 -Day 21
$ws.Range("E16").NumberFormat = "mm-dd-yy"
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 2

# Extend the total-hours formula to include the new row
$ws.Range("F17").Formula = "=SUM(F5:F16)"

# Update the view: scroll so row 7 is the top-left, and select H16
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H16").Select()
